# Roller shades rule refactor
# Remove the "G" roller-shade entry (row 9: G_Motion_Detector /
# G_Brightness_Sensor / G_Lights_Analog) from the motion-template-data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9:C9").ClearContents()
